$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Nodes")

# Update H57:H228 from 0.6 to 0.5, and I57:I228 from 0.4 to 0.5
for ($r = 57; $r -le 228; $r++) {
    $ws.Cells.Item($r, 8).Value = 0.5
    $ws.Cells.Item($r, 9).Value = 0.5
}

# Update the sheet view: frozen pane top-left cell and selection
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 46
$ws.Range("I57:I309").Select()
$excel.ActiveWindow.RangeSelection.Item(1,1)
$excel.ActiveCell
